$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-7 (columns F:AO) with new odds ---
# Row 2
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 1.73
$ws.Range("I2").Value = 1.78
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 3.9
$ws.Range("L2").Value = 1.49
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 2.96
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 1.66
$ws.Range("Q2").Value = 2.26
$ws.Range("R2").Value = 1.23
$ws.Range("S2").Value = 4.5
$ws.Range("T2").Value = 2.16
$ws.Range("U2").Value = 1.72
$ws.Range("V2").Value = 2.04
$ws.Range("W2").Value = 1.15
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 6.8
$ws.Range("Z2").Value = 9.4
$ws.Range("AA2").Value = 19
$ws.Range("AB2").Value = 20
$ws.Range("AC2").Value = 8.800000000000001
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 23
$ws.Range("AF2").Value = 55
$ws.Range("AG2").Value = 24
$ws.Range("AH2").Value = 28
$ws.Range("AI2").Value = 55
$ws.Range("AJ2").Value = 230
$ws.Range("AK2").Value = 130
$ws.Range("AL2").Value = 140
$ws.Range("AM2").Value = 250
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 16.5

# Row 3
$ws.Range("F3").Value = 7.6
$ws.Range("G3").Value = 9.199999999999999
$ws.Range("H3").Value = 1.54
$ws.Range("I3").Value = 1.59
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 4.5
$ws.Range("L3").Value = 1.45
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 3.15
$ws.Range("O3").Value = 1.39
$ws.Range("P3").Value = 1.74
$ws.Range("Q3").Value = 1.97
$ws.Range("R3").Value = 1.28
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 2.22
$ws.Range("U3").Value = 1.71
$ws.Range("V3").Value = 2.34
$ws.Range("W3").Value = 1.1
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 7.8
$ws.Range("Z3").Value = 9.800000000000001
$ws.Range("AA3").Value = 17.5
$ws.Range("AB3").Value = 25
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 12.5
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 980
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 200
$ws.Range("AL3").Value = 180
$ws.Range("AM3").Value = 270
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 14

# Row 4
$ws.Range("F4").Value = 3.3
$ws.Range("G4").Value = 4.1
$ws.Range("H4").Value = 2.22
$ws.Range("I4").Value = 2.5
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 3.65
$ws.Range("L4").Value = 1.45
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 3.1
$ws.Range("O4").Value = 1.39
$ws.Range("P4").Value = 1.73
$ws.Range("Q4").Value = 2.12
$ws.Range("R4").Value = 1.28
$ws.Range("S4").Value = 3.9
$ws.Range("T4").Value = 1.84
$ws.Range("U4").Value = 1.96
$ws.Range("V4").Value = 1.66
$ws.Range("W4").Value = 1.34
$ws.Range("X4").Value = 14.5
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 14.5
$ws.Range("AC4").Value = 9.199999999999999
$ws.Range("AD4").Value = 13.5
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 60
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 140
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 26

# Row 5
$ws.Range("F5").Value = 1.63
$ws.Range("G5").Value = 1.69
$ws.Range("H5").Value = 7.2
$ws.Range("I5").Value = 7.4
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 1.5
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 3.15
$ws.Range("O5").Value = 1.39
$ws.Range("P5").Value = 1.73
$ws.Range("Q5").Value = 2.18
$ws.Range("R5").Value = 1.26
$ws.Range("S5").Value = 4
$ws.Range("T5").Value = 2.12
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 1.15
$ws.Range("W5").Value = 2.4
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 19.5
$ws.Range("Z5").Value = 55
$ws.Range("AA5").Value = 240
$ws.Range("AB5").Value = 7
$ws.Range("AC5").Value = 8.800000000000001
$ws.Range("AD5").Value = 28
$ws.Range("AE5").Value = 130
$ws.Range("AF5").Value = 8.800000000000001
$ws.Range("AG5").Value = 10.5
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 140
$ws.Range("AJ5").Value = 16
$ws.Range("AK5").Value = 20
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 210
$ws.Range("AN5").Value = 12.5
$ws.Range("AO5").Value = 210

# Row 6
$ws.Range("F6").Value = 2.28
$ws.Range("G6").Value = 2.48
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 4.3
$ws.Range("J6").Value = 2.92
$ws.Range("K6").Value = 3.2
$ws.Range("L6").Value = 1.63
$ws.Range("M6").Value = 1.12
$ws.Range("N6").Value = 2.38
$ws.Range("O6").Value = 1.61
$ws.Range("P6").Value = 1.45
$ws.Range("Q6").Value = 2.8
$ws.Range("R6").Value = 1.17
$ws.Range("S6").Value = 6
$ws.Range("T6").Value = 2.24
$ws.Range("U6").Value = 1.69
$ws.Range("V6").Value = 1.31
$ws.Range("W6").Value = 1.68
$ws.Range("X6").Value = 9.199999999999999
$ws.Range("Y6").Value = 980
$ws.Range("Z6").Value = 34
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 8
$ws.Range("AC6").Value = 8.6
$ws.Range("AD6").Value = 22
$ws.Range("AE6").Value = 90
$ws.Range("AF6").Value = 16
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 32
$ws.Range("AI6").Value = 130
$ws.Range("AJ6").Value = 44
$ws.Range("AK6").Value = 46
$ws.Range("AL6").Value = 90
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("F7").Value = 1.87
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 3.15
$ws.Range("K7").Value = 3.5
$ws.Range("L7").Value = 1.56
$ws.Range("M7").Value = 1.12
$ws.Range("N7").Value = 2.6
$ws.Range("O7").Value = 1.54
$ws.Range("P7").Value = 1.51
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.18
$ws.Range("S7").Value = 5.3
$ws.Range("T7").Value = 2.28
$ws.Range("U7").Value = 1.64
$ws.Range("V7").Value = 1.2
$ws.Range("W7").Value = 2
$ws.Range("X7").Value = 8.800000000000001
$ws.Range("Y7").Value = 14.5
$ws.Range("Z7").Value = 42
$ws.Range("AA7").Value = 200
$ws.Range("AB7").Value = 6.8
$ws.Range("AC7").Value = 8.199999999999999
$ws.Range("AD7").Value = 25
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 10.5
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 30
$ws.Range("AI7").Value = 160
$ws.Range("AJ7").Value = 24
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 70
$ws.Range("AM7").Value = 300
$ws.Range("AN7").Value = 23
$ws.Range("AO7").Value = 1000

# --- Add new rows 8-11 ---
# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "Brazilian Serie A"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-10-08"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "21:00:00"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "Mirassol"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "Fluminense"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = 2.36
$ws.Range("G8").Value = 2.56
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 3.05
$ws.Range("K8").Value = 3.3
$ws.Range("L8").Value = 1.56
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 2.82
$ws.Range("O8").Value = 1.45
$ws.Range("P8").Value = 1.63
$ws.Range("Q8").Value = 2.38
$ws.Range("R8").Value = 1.22
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = 1.96
$ws.Range("U8").Value = 1.86
$ws.Range("V8").Value = 1.35
$ws.Range("W8").Value = 1.64
$ws.Range("X8").Value = 10
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 29
$ws.Range("AA8").Value = 80
$ws.Range("AB8").Value = 10
$ws.Range("AC8").Value = 8.4
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 55
$ws.Range("AF8").Value = 17
$ws.Range("AG8").Value = 14.5
$ws.Range("AH8").Value = 26
$ws.Range("AI8").Value = 85
$ws.Range("AJ8").Value = 42
$ws.Range("AK8").Value = 38
$ws.Range("AL8").Value = 60
$ws.Range("AM8").Value = 190
$ws.Range("AN8").Value = 36
$ws.Range("AO8").Value = 70

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "Brazilian Serie B"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2025-10-08"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "21:00:00"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "Avai"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "Volta Redonda"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = 1.7
$ws.Range("G9").Value = 1.8
$ws.Range("H9").Value = 5.4
$ws.Range("I9").Value = 7.4
$ws.Range("J9").Value = 3.65
$ws.Range("K9").Value = 4.2
$ws.Range("L9").Value = 1.43
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 3.3
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 1.79
$ws.Range("Q9").Value = 2.06
$ws.Range("R9").Value = 1.29
$ws.Range("S9").Value = 3.7
$ws.Range("T9").Value = 1.98
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.17
$ws.Range("W9").Value = 2.24
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 19
$ws.Range("Z9").Value = 48
$ws.Range("AA9").Value = 190
$ws.Range("AB9").Value = 8
$ws.Range("AC9").Value = 9.199999999999999
$ws.Range("AD9").Value = 25
$ws.Range("AE9").Value = 110
$ws.Range("AF9").Value = 10.5
$ws.Range("AG9").Value = 11
$ws.Range("AH9").Value = 25
$ws.Range("AI9").Value = 120
$ws.Range("AJ9").Value = 19
$ws.Range("AK9").Value = 21
$ws.Range("AL9").Value = 44
$ws.Range("AM9").Value = 190
$ws.Range("AN9").Value = 13.5
$ws.Range("AO9").Value = 150

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "Brazilian Serie B"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "2025-10-08"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "21:30:00"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "America MG"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "Vila Nova"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = 1.98
$ws.Range("G10").Value = 2.12
$ws.Range("H10").Value = 4.5
$ws.Range("I10").Value = 5.1
$ws.Range("J10").Value = 3.15
$ws.Range("K10").Value = 3.45
$ws.Range("L10").Value = 1.59
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 2.68
$ws.Range("O10").Value = 1.52
$ws.Range("P10").Value = 1.56
$ws.Range("Q10").Value = 2.52
$ws.Range("R10").Value = 1.2
$ws.Range("S10").Value = 5.1
$ws.Range("T10").Value = 2.18
$ws.Range("U10").Value = 1.73
$ws.Range("V10").Value = 1.24
$ws.Range("W10").Value = 1.89
$ws.Range("X10").Value = 10.5
$ws.Range("Y10").Value = 13
$ws.Range("Z10").Value = 34
$ws.Range("AA10").Value = 150
$ws.Range("AB10").Value = 7
$ws.Range("AC10").Value = 7.8
$ws.Range("AD10").Value = 21
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 11
$ws.Range("AG10").Value = 11.5
$ws.Range("AH10").Value = 26
$ws.Range("AI10").Value = 130
$ws.Range("AJ10").Value = 25
$ws.Range("AK10").Value = 28
$ws.Range("AL10").Value = 60
$ws.Range("AM10").Value = 260
$ws.Range("AN10").Value = 28
$ws.Range("AO10").Value = 1000

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "US MLS"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2025-10-08"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "23:30:00"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "Los Angeles FC"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "Toronto FC"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = 1.56
$ws.Range("G11").Value = 1.61
$ws.Range("H11").Value = 6
$ws.Range("I11").Value = 6.8
$ws.Range("J11").Value = 4.6
$ws.Range("K11").Value = 5
$ws.Range("L11").Value = 1.31
$ws.Range("M11").Value = 1.03
$ws.Range("N11").Value = 5.1
$ws.Range("O11").Value = 1.21
$ws.Range("P11").Value = 2.38
$ws.Range("Q11").Value = 1.63
$ws.Range("R11").Value = 1.55
$ws.Range("S11").Value = 2.58
$ws.Range("T11").Value = 1.75
$ws.Range("U11").Value = 2.18
$ws.Range("V11").Value = 1.17
$ws.Range("W11").Value = 2.62
$ws.Range("X11").Value = 22
$ws.Range("Y11").Value = 30
$ws.Range("Z11").Value = 55
$ws.Range("AA11").Value = 200
$ws.Range("AB11").Value = 11
$ws.Range("AC11").Value = 11.5
$ws.Range("AD11").Value = 24
$ws.Range("AE11").Value = 80
$ws.Range("AF11").Value = 11
$ws.Range("AG11").Value = 10
$ws.Range("AH11").Value = 20
$ws.Range("AI11").Value = 75
$ws.Range("AJ11").Value = 16
$ws.Range("AK11").Value = 16
$ws.Range("AL11").Value = 30
$ws.Range("AM11").Value = 95
$ws.Range("AN11").Value = 7.2
$ws.Range("AO11").Value = 75
